$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1002743.6
$ws.Range("J17").Value = 1002743.6
$ws.Range("L17").Value = 3008230.8
$ws.Range("N17").Value = -3008566.8
$ws.Range("H19").Value = 651
$ws.Range("I19").Value = 499.0909
$ws.Range("J19").Value = 802.9091
$ws.Range("K19").Value = 499.0909
$ws.Range("L19").Value = 802.9091
$ws.Range("M19").Value = -324.0909
$ws.Range("N19").Value = -1152.9091
$ws.Range("H20").Value = 34007
$ws.Range("I20").Value = 34007
$ws.Range("K20").Value = 34007
$ws.Range("M20").Value = -33777
$ws.Range("H21").Value = 14725
$ws.Range("J21").Value = 27500
$ws.Range("L21").Value = 27500
$ws.Range("N21").Value = -28436
$ws.Range("H23").Value = 14725
$ws.Range("J23").Value = 27500
$ws.Range("L23").Value = 27500
$ws.Range("N23").Value = -27968
$ws.Range("H26").Value = 33000
$ws.Range("J26").Value = 48000
$ws.Range("L26").Value = 48000
$ws.Range("N26").Value = -48688
$ws.Range("H35").Value = 34007
$ws.Range("I35").Value = 34007
$ws.Range("K35").Value = 34007
$ws.Range("M35").Value = -33628
$ws.Range("H40").Value = 2265.4666
$ws.Range("I40").Value = 2398.75
$ws.Range("J40").Value = 2113.1428
$ws.Range("K40").Value = 2398.75
$ws.Range("L40").Value = 2113.1428
$ws.Range("M40").Value = -2223.75
$ws.Range("N40").Value = -2463.1428
$ws.Range("H76").Value = 5294476
$ws.Range("I76").Value = 11114294
$ws.Range("J76").Value = 3731.6365
$ws.Range("K76").Value = 11114294
$ws.Range("L76").Value = 3731.6365
$ws.Range("M76").Value = -11113979
$ws.Range("N76").Value = -4361.636500000001
$ws.Range("H79").Value = 5294476
$ws.Range("I79").Value = 11114294
$ws.Range("J79").Value = 3731.6365
$ws.Range("K79").Value = 11114294
$ws.Range("L79").Value = 3731.6365
$ws.Range("M79").Value = -11113202
$ws.Range("N79").Value = -5915.636500000001
$ws.Range("H112").Value = 24794790
$ws.Range("I112").Value = 993
$ws.Range("J112").Value = 34092464
$ws.Range("K112").Value = 2979
$ws.Range("L112").Value = 102277392
$ws.Range("M112").Value = -1871
$ws.Range("N112").Value = -102279608
$ws.Range("H133").Value = 42596.125
$ws.Range("J133").Value = 42596.125
$ws.Range("L133").Value = 42596.125
$ws.Range("N133").Value = -52716.125
$ws.Range("H135").Value = 1137.3715
$ws.Range("I135").Value = 1051.871
$ws.Range("J135").Value = 1800
$ws.Range("K135").Value = 9466.839
$ws.Range("L135").Value = 16200
$ws.Range("M135").Value = -6931.839
$ws.Range("N135").Value = -21270
$ws.Range("H137").Value = 33334352
$ws.Range("I137").Value = 45455500
$ws.Range("K137").Value = 136366500
$ws.Range("M137").Value = -136363950

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 26500
$ws.Range("J12").Value = 3000
$ws.Range("L12").Value = 3000
$ws.Range("N12").Value = -3346
$ws.Range("H19").Value = 30500
$ws.Range("I19").Value = 1000
$ws.Range("J19").Value = 60000
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 60000
$ws.Range("M19").Value = -771
$ws.Range("N19").Value = -60458
$ws.Range("H32").Value = 25762.268
$ws.Range("I32").Value = 3125.9048
$ws.Range("J32").Value = 342671.34
$ws.Range("K32").Value = 3125.9048
$ws.Range("L32").Value = 342671.34
$ws.Range("M32").Value = -2838.9048
$ws.Range("N32").Value = -343245.34
$ws.Range("H122").Value = 1714.6364
$ws.Range("I122").Value = 1646.9615
$ws.Range("J122").Value = 1966
$ws.Range("K122").Value = 4940.8845
$ws.Range("L122").Value = 5898
$ws.Range("M122").Value = -2490.8845
$ws.Range("N122").Value = -10798
$ws.Range("H123").Value = 45000
$ws.Range("J123").Value = 45000
$ws.Range("L123").Value = 45000
$ws.Range("N123").Value = -54800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 17316.857
$ws.Range("I86").Value = 3535.1667
$ws.Range("J86").Value = 100007
$ws.Range("K86").Value = 3535.1667
$ws.Range("L86").Value = 100007
$ws.Range("M86").Value = -2412.1667
$ws.Range("N86").Value = -102253
$ws.Range("H89").Value = 17316.857
$ws.Range("I89").Value = 3535.1667
$ws.Range("J89").Value = 100007
$ws.Range("K89").Value = 17675.8335
$ws.Range("L89").Value = 500035
$ws.Range("M89").Value = -12059.8335
$ws.Range("N89").Value = -511267
$ws.Range("H94").Value = 1860.1818
$ws.Range("I94").Value = 1164.6666
$ws.Range("J94").Value = 4990
$ws.Range("K94").Value = 1164.6666
$ws.Range("L94").Value = 4990
$ws.Range("M94").Value = -713.6666
$ws.Range("N94").Value = -5892
$ws.Range("H107").Value = 1594.6875
$ws.Range("I107").Value = 1374.0834
$ws.Range("J107").Value = 2256.5
$ws.Range("K107").Value = 1374.0834
$ws.Range("L107").Value = 2256.5
$ws.Range("M107").Value = 545.9166
$ws.Range("N107").Value = -6096.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1169.1666
$ws.Range("I19").Value = 403
$ws.Range("J19").Value = 5000
$ws.Range("K19").Value = 403
$ws.Range("L19").Value = 5000
$ws.Range("M19").Value = -233
$ws.Range("N19").Value = -5340
$ws.Range("H24").Value = 1169.1666
$ws.Range("I24").Value = 403
$ws.Range("J24").Value = 5000
$ws.Range("K24").Value = 403
$ws.Range("L24").Value = 5000
$ws.Range("M24").Value = -233
$ws.Range("N24").Value = -5340
$ws.Range("H31").Value = 4036.465
$ws.Range("I31").Value = 1065.4412
$ws.Range("J31").Value = 15260.333
$ws.Range("K31").Value = 1065.4412
$ws.Range("L31").Value = 15260.333
$ws.Range("M31").Value = -770.4412
$ws.Range("N31").Value = -15850.333
$ws.Range("H34").Value = 4036.465
$ws.Range("I34").Value = 1065.4412
$ws.Range("J34").Value = 15260.333
$ws.Range("K34").Value = 1065.4412
$ws.Range("L34").Value = 15260.333
$ws.Range("M34").Value = -863.4412
$ws.Range("N34").Value = -15664.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 57500
$ws.Range("J37").Value = 57500
$ws.Range("L37").Value = 172500
$ws.Range("N37").Value = -172724
$ws.Range("H56").Value = 5422.5
$ws.Range("I56").Value = 5422.5
$ws.Range("K56").Value = 5422.5
$ws.Range("M56").Value = -4892.5
$ws.Range("H81").Value = 11135.454
$ws.Range("I81").Value = 1483.3334
$ws.Range("J81").Value = 22718
$ws.Range("K81").Value = 4450.0002
$ws.Range("L81").Value = 68154
$ws.Range("M81").Value = -3327.0002
$ws.Range("N81").Value = -70400
$ws.Range("H84").Value = 11135.454
$ws.Range("I84").Value = 1483.3334
$ws.Range("J84").Value = 22718
$ws.Range("K84").Value = 13350.0006
$ws.Range("L84").Value = 204462
$ws.Range("M84").Value = -7734.000599999999
$ws.Range("N84").Value = -215694

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 8000
$ws.Range("J43").Value = 8000
$ws.Range("L43").Value = 8000
$ws.Range("N43").Value = -8302
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 1000
$ws.Range("K97").Value = 1000
$ws.Range("M97").Value = -504
$ws.Range("H123").Value = 9884.105
$ws.Range("J123").Value = 9884.105
$ws.Range("L123").Value = 9884.105
$ws.Range("N123").Value = -14784.105

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 13101.5
$ws.Range("I19").Value = 17235.334
$ws.Range("J19").Value = 700
$ws.Range("K19").Value = 17235.334
$ws.Range("L19").Value = 700
$ws.Range("M19").Value = -17065.334
$ws.Range("N19").Value = -1040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 8185
$ws.Range("I7").Value = 9215.833000000001
$ws.Range("J7").Value = 2000
$ws.Range("K7").Value = 9215.833000000001
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = -9102.833000000001
$ws.Range("N7").Value = -2226
$ws.Range("H132").Value = 2332.4146
$ws.Range("I132").Value = 1934.5385
$ws.Range("J132").Value = 3022.0667
$ws.Range("K132").Value = 5803.6155
$ws.Range("L132").Value = 9066.2001
$ws.Range("M132").Value = -3273.6155
$ws.Range("N132").Value = -14126.2001
$ws.Range("H136").Value = 2010.1765
$ws.Range("I136").Value = 955.2727
$ws.Range("J136").Value = 3944.1667
$ws.Range("K136").Value = 2865.8181
$ws.Range("L136").Value = 11832.5001
$ws.Range("M136").Value = -315.8181
$ws.Range("N136").Value = -16932.5001
